# Update TPM-derived values in the LR-pairs sheet (Wnt7b-Lrp5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05003266666666667
$ws.Range("H2").Value = 0.150098
$ws.Range("M2").Value = 18.42392
$ws.Range("N2").Value = 55.27176
$ws.Range("O2").Value = 0.3903243738016154
$ws.Range("P2").Value = 0.3903243738016154
$ws.Range("Q2").Value = 0.9217978480533333
$ws.Range("R2").Value = 8.29618063248
$ws.Range("S2").Value = 0.3903243738016154
$ws.Range("T2").Value = 0.3903243738016154

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05003266666666667
$ws.Range("H3").Value = 0.150098
$ws.Range("O3").Value = 0.2625687066780312
$ws.Range("P3").Value = 0.2625687066780312
$ws.Range("Q3").Value = 0.6200875093313334
$ws.Range("R3").Value = 5.580787583982
$ws.Range("S3").Value = 0.2625687066780312
$ws.Range("T3").Value = 0.2625687066780312

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05003266666666667
$ws.Range("H4").Value = 0.150098
$ws.Range("M4").Value = 10.76369066666667
$ws.Range("N4").Value = 32.291072
$ws.Range("O4").Value = 0.2280367489253622
$ws.Range("P4").Value = 0.2280367489253622
$ws.Range("Q4").Value = 0.5385361472284445
$ws.Range("R4").Value = 4.846825325056
$ws.Range("S4").Value = 0.2280367489253622
$ws.Range("T4").Value = 0.2280367489253622

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05003266666666667
$ws.Range("H5").Value = 0.150098
$ws.Range("M5").Value = 5.620297999999999
$ws.Range("N5").Value = 16.860894
$ws.Range("O5").Value = 0.1190701705949913
$ws.Range("P5").Value = 0.1190701705949913
$ws.Range("Q5").Value = 0.2811984964013333
$ws.Range("R5").Value = 2.530786467612
$ws.Range("S5").Value = 0.1190701705949913
$ws.Range("T5").Value = 0.1190701705949913
